$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.505.16'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '3.598.89'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '609.46'
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("D6").Value = '148.98'
$ws.Range("E6").Value = '  +3.25%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D12").Value = '4.212.41'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("D14").Value = '29.84'
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").Value = '3.608.57'
$ws.Range("E15").Value = '  +1.61%  '
$ws.Range("D16").Value = '66.606.64'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("E18").Value = '  +2.10%  '
$ws.Range("D19").Value = '6.39'
$ws.Range("D20").Value = '15.12'
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("D21").Value = '428.00'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").Value = '78.80'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '0.0000123'
$ws.Range("E25").Value = '  +4.47%  '
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").Value = '8.33'
$ws.Range("E26").Value = '  +4.90%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '9.43'
$ws.Range("E27").Value = '  +4.01%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '2.52'
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = '1.48'
$ws.Range("E30").Value = '  +1.16%  '
$ws.Range("B31").Value = 'RenzoRestakedETH'
$ws.Range("C31").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D31").Value = '3.596.63'
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '25.48'
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").Value = '0.157'
$ws.Range("E33").Value = '  +4.40%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '7.87'
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '5.68'
$ws.Range("E36").Value = '  +1.40%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '1.70'
$ws.Range("E37").Value = '  -2.01%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = '177.29'
$ws.Range("E38").Value = '  +1.23%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.0857'
$ws.Range("E39").Value = '  +1.00%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").Value = '5.24'
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").Value = '0.900'
$ws.Range("E41").Value = '  +0.86%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '1.91'
$ws.Range("E42").Value = '  -0.67%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '2.58'
$ws.Range("E43").Value = '  +10.32%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").Value = '1.18'
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '25.06'
$ws.Range("E46").Value = '  -2.57%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '24.07'
$ws.Range("E47").Value = '  +1.99%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '7.21'
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("B49").Value = 'SuiNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D49").Value = '0.955'
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.429.27'
$ws.Range("E50").Value = '  +5.67%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '0.236'
$ws.Range("E51").Value = '  -0.30%  '
